$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new Fitness (column C) values per the corrected SA algorithm run.
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 3).Value = 7581
}

for ($r = 12; $r -le 217; $r++) {
    $ws.Cells.Item($r, 3).Value = 7534
}

for ($r = 218; $r -le 252; $r++) {
    $ws.Cells.Item($r, 3).Value = 7320
}
